$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new rows 16-18 with date (col A) and time (col B) values,
# matching the existing formatting used by rows 2-15 (style carries
# over automatically since these cells inherit the column style).
$ws.Range("A16").Value = 45201
$ws.Range("B16").Value = 0.40625

$ws.Range("A17").Value = 45202
$ws.Range("B17").Value = 0.40277777777777773

$ws.Range("A18").Value = 45203
$ws.Range("B18").Value = 0.40625

# Match the styles of row 15 (date format in A, time format in B) by
# copying the cell formatting rather than re-deriving a number format
# string, so the original style indexes (numFmtId 14 / 20) are reused
# instead of creating new custom formats.
$ws.Range("A15:B15").Copy() | Out-Null
$ws.Range("A16:B18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Move the selection to reflect where the user ended up after entry
$ws.Range("B19").Select()
